$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.191.79"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "1.906.55"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5228"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3768"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07270"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9069"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08471"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.87%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.918.23"
$ws.Range("E13").Value = "  +1.87%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.302"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008680"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").Value = "27.231.26"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.103"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").Value = "2.158.94"
$ws.Range("E22").Value = "  +1.41%  "

$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.448"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.334"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.765"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.832"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.927"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09320"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.28%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05061"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.249"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("E36").Value = "  +4.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.943"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.597"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5741"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02008"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.112"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.622"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "115.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1520"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("E46").Value = "  +1.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.628"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("E51").Value = "  -0.05%  "
